# Auto-generated Excel COM-interop script
# Applies numeric value updates per the commit diff (currentAveragePrice /
# LevePrice / LeveProfit columns refreshed by the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H62").Value = 4664.0835
$ws_ALC.Range("I62").Value = 4842.727
$ws_ALC.Range("K62").Value = 4842.727
$ws_ALC.Range("M62").Value = -4218.727

$ws_ALC.Range("H65").Value = 4664.0835
$ws_ALC.Range("I65").Value = 4842.727
$ws_ALC.Range("K65").Value = 24213.635
$ws_ALC.Range("M65").Value = -21093.635

$ws_ALC.Range("H98").Value = 1353.762
$ws_ALC.Range("I98").Value = 1259.1052
$ws_ALC.Range("J98").Value = 2253
$ws_ALC.Range("K98").Value = 1259.1052
$ws_ALC.Range("L98").Value = 2253
$ws_ALC.Range("M98").Value = 238.8948
$ws_ALC.Range("N98").Value = -5249

$ws_ALC.Range("H106").Value = 19976
$ws_ALC.Range("I106").Value = 20732.688
$ws_ALC.Range("K106").Value = 20732.688
$ws_ALC.Range("M106").Value = -20101.688

$ws_ALC.Range("H112").Value = 8879.727999999999
$ws_ALC.Range("J112").Value = 8879.727999999999
$ws_ALC.Range("L112").Value = 26639.184
$ws_ALC.Range("N112").Value = -28855.184

$ws_ALC.Range("H122").Value = 1353.762
$ws_ALC.Range("I122").Value = 1259.1052
$ws_ALC.Range("J122").Value = 2253
$ws_ALC.Range("K122").Value = 3777.3156
$ws_ALC.Range("L122").Value = 6759
$ws_ALC.Range("M122").Value = -1327.3156
$ws_ALC.Range("N122").Value = -11659

$ws_ALC.Range("H127").Value = 584.3125
$ws_ALC.Range("I127").Value = 475.46667
$ws_ALC.Range("K127").Value = 1426.40001
$ws_ALC.Range("M127").Value = 3533.59999

$ws_ALC.Range("H132").Value = 3788.7273
$ws_ALC.Range("I132").Value = 3417.4614
$ws_ALC.Range("K132").Value = 10252.3842
$ws_ALC.Range("M132").Value = -7722.3842

$ws_ALC.Range("H138").Value = 3956.9348
$ws_ALC.Range("J138").Value = 4579.36
$ws_ALC.Range("L138").Value = 13738.08
$ws_ALC.Range("N138").Value = -24018.08

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H94").Value = 74232.586
$ws_ARM.Range("J94").Value = 74559.56
$ws_ARM.Range("L94").Value = 74559.56
$ws_ARM.Range("N94").Value = -76361.56

$ws_ARM.Range("H102").Value = 3023
$ws_ARM.Range("I102").Value = 3317.4546
$ws_ARM.Range("K102").Value = 3317.4546
$ws_ARM.Range("M102").Value = -1695.4546

$ws_ARM.Range("H141").Value = 70047.664
$ws_ARM.Range("J141").Value = 70047.664
$ws_ARM.Range("L141").Value = 70047.664
$ws_ARM.Range("N141").Value = -80407.664

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H99").Value = 3073.6
$ws_BSM.Range("I99").Value = 3118.7568
$ws_BSM.Range("K99").Value = 3118.7568
$ws_BSM.Range("M99").Value = -1620.7568

$ws_BSM.Range("H134").Value = 1983.762
$ws_BSM.Range("I134").Value = 1987.1666
$ws_BSM.Range("K134").Value = 5961.4998
$ws_BSM.Range("M134").Value = -3426.4998

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H28").Value = 13695.818
$ws_CRP.Range("J28").Value = 13695.818
$ws_CRP.Range("L28").Value = 13695.818
$ws_CRP.Range("N28").Value = -14185.818

$ws_CRP.Range("H31").Value = 474510.94
$ws_CRP.Range("I31").Value = 14552
$ws_CRP.Range("J31").Value = 641768.75
$ws_CRP.Range("K31").Value = 14552
$ws_CRP.Range("L31").Value = 641768.75
$ws_CRP.Range("M31").Value = -14257
$ws_CRP.Range("N31").Value = -642358.75

$ws_CRP.Range("H34").Value = 474510.94
$ws_CRP.Range("I34").Value = 14552
$ws_CRP.Range("J34").Value = 641768.75
$ws_CRP.Range("K34").Value = 14552
$ws_CRP.Range("L34").Value = 641768.75
$ws_CRP.Range("M34").Value = -14350
$ws_CRP.Range("N34").Value = -642172.75

$ws_CRP.Range("H132").Value = 1875
$ws_CRP.Range("I132").Value = 1387.3334
$ws_CRP.Range("K132").Value = 4162.0002
$ws_CRP.Range("M132").Value = -1632.0002

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H59").Value = 3547.25
$ws_CUL.Range("I59").Value = 4563
$ws_CUL.Range("K59").Value = 13689
$ws_CUL.Range("M59").Value = -13149

$ws_CUL.Range("H103").Value = 890.8570999999999
$ws_CUL.Range("J103").Value = 1499.5
$ws_CUL.Range("L103").Value = 4498.5
$ws_CUL.Range("N103").Value = -6256.5

$ws_CUL.Range("H113").Value = 8548692
$ws_CUL.Range("I113").Value = 1952.7858
$ws_CUL.Range("J113").Value = 13334867
$ws_CUL.Range("K113").Value = 5858.357400000001
$ws_CUL.Range("L113").Value = 40004601
$ws_CUL.Range("M113").Value = -3688.357400000001
$ws_CUL.Range("N113").Value = -40008941

$ws_CUL.Range("H132").Value = 3988.8462
$ws_CUL.Range("I132").Value = 3612.9092
$ws_CUL.Range("J132").Value = 4264.533
$ws_CUL.Range("K132").Value = 32516.1828
$ws_CUL.Range("L132").Value = 38380.79700000001
$ws_CUL.Range("M132").Value = -29986.1828
$ws_CUL.Range("N132").Value = -43440.79700000001

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H99").Value = 10499
$ws_GSM.Range("I99").Value = 10499
$ws_GSM.Range("J99").Value = 0
$ws_GSM.Range("K99").Value = 10499
$ws_GSM.Range("L99").Value = 0
$ws_GSM.Range("M99").Value = -8253
$ws_GSM.Range("N99").ClearContents()

$ws_GSM.Range("H105").Value = 40000
$ws_GSM.Range("J105").Value = 40000
$ws_GSM.Range("L105").Value = 40000
$ws_GSM.Range("N105").Value = -46988

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H22").Value = 41356.52
$ws_LTW.Range("I22").Value = 125427.5
$ws_LTW.Range("J22").Value = 1793.7059
$ws_LTW.Range("K22").Value = 125427.5
$ws_LTW.Range("L22").Value = 1793.7059
$ws_LTW.Range("M22").Value = -125132.5
$ws_LTW.Range("N22").Value = -2383.7059

$ws_LTW.Range("H27").Value = 41356.52
$ws_LTW.Range("I27").Value = 125427.5
$ws_LTW.Range("J27").Value = 1793.7059
$ws_LTW.Range("K27").Value = 125427.5
$ws_LTW.Range("L27").Value = 1793.7059
$ws_LTW.Range("M27").Value = -125320.5
$ws_LTW.Range("N27").Value = -2007.7059

$ws_LTW.Range("H55").Value = 2044.6316
$ws_LTW.Range("I55").Value = 986.4545000000001
$ws_LTW.Range("J55").Value = 3499.625
$ws_LTW.Range("K55").Value = 986.4545000000001
$ws_LTW.Range("L55").Value = 3499.625
$ws_LTW.Range("M55").Value = -813.4545000000001
$ws_LTW.Range("N55").Value = -3845.625

$ws_LTW.Range("H106").Value = 39788
$ws_LTW.Range("J106").Value = 39788
$ws_LTW.Range("L106").Value = 39788
$ws_LTW.Range("N106").Value = -42312

$ws_LTW.Range("H132").Value = 5537.8066
$ws_LTW.Range("I132").Value = 5265.778
$ws_LTW.Range("J132").Value = 7374
$ws_LTW.Range("K132").Value = 15797.334
$ws_LTW.Range("L132").Value = 22122
$ws_LTW.Range("M132").Value = -13267.334
$ws_LTW.Range("N132").Value = -27182

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H54").Value = 49999.91

$ws_WVR.Range("H104").Value = 71250
$ws_WVR.Range("J104").Value = 71250
$ws_WVR.Range("L104").Value = 71250
$ws_WVR.Range("N104").Value = -78238

$ws_WVR.Range("H122").Value = 16381979
$ws_WVR.Range("I122").Value = 17047906
$ws_WVR.Range("K122").Value = 51143718
$ws_WVR.Range("M122").Value = -51141268

$ws_WVR.Range("H132").Value = 2014.125
$ws_WVR.Range("I132").Value = 1900.1063
$ws_WVR.Range("J132").Value = 2609.5557
$ws_WVR.Range("K132").Value = 5700.3189
$ws_WVR.Range("L132").Value = 7828.6671
$ws_WVR.Range("M132").Value = -3170.3189
$ws_WVR.Range("N132").Value = -12888.6671

$ws_WVR.Range("H136").Value = 4554.4844
$ws_WVR.Range("I136").Value = 4370.5815
$ws_WVR.Range("K136").Value = 13111.7445
$ws_WVR.Range("M136").Value = -10561.7445
